# Add the new "ReservasServicios" worksheet as the last sheet in the workbook
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ReservasServicios"

# Header row (row 1)
$ws.Range("A1").Value = "ID Reserva"
$ws.Range("B1").Value = "Cliente"
$ws.Range("C1").Value = "Servicio"
$ws.Range("D1").Value = "Fecha Reserva"
$ws.Range("E1").Value = "Fecha Servicio"

# Data row (row 2)
$ws.Range("A2").Value = 12
$ws.Range("B2").Value = "mgiue"
$ws.Range("C2").Value = 345
$ws.Range("D2").Value = "28/05/2024 14:25"
$ws.Range("E2").Value = "28/05/2024 14:30"

# Column widths (account for the internal +5/6 character padding Excel applies
# on top of the ColumnWidth property so the stored OOXML width matches exactly)
$ws.Columns.Item(1).ColumnWidth = 12 - 5/6
$ws.Columns.Item(2).ColumnWidth = 9 - 5/6
$ws.Columns.Item(3).ColumnWidth = 10 - 5/6
$ws.Columns.Item(4).ColumnWidth = 18 - 5/6
$ws.Columns.Item(5).ColumnWidth = 18 - 5/6

# Apply the same bold/centered/bordered header style used by the other sheets
# (e.g. Inventario!A1) to the header row, reusing the existing cell style.
$srcSheet = $wb.Worksheets.Item("Inventario")
$srcSheet.Range("A1").Copy()
$ws.Range("A1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A1").Select() | Out-Null
